$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected financial figures for actual fiscal years (rows 2-6, FY2014-FY2018)
$rowData = @{}
$rowData[2] = @{ "D"=2619; "E"=54; "F"=54; "G"=48; "H"=50; "I"=54; "J"=-3; "K"=1256; "L"=715; "M"=541; "N"=533; "O"=7; "P"=300; "Q"=-9; "R"=-14; "S"=39; "T"=26; "U"=-36; "V"=149; "W"=2.07; "X"=1.93; "Y"=11.16; "Z"=4.33; "AA"=132.25; "AB"=95.58; "AC"=89; "AD"=16.98; "AE"=964; "AF"=1.57; "AG"=20; "AH"=1.32; "AI"=20.66; "AJ"=59991641 }
$rowData[3] = @{ "D"=3083; "E"=81; "F"=81; "G"=82; "H"=85; "I"=85; "J"=0; "K"=1263; "L"=658; "M"=605; "N"=597; "O"=7; "P"=300; "Q"=217; "R"=-43; "S"=-55; "T"=50; "U"=167; "V"=105; "W"=2.63; "X"=2.74; "Y"=15.04; "Z"=6.72; "AA"=108.87; "AB"=117.08; "AC"=142; "AD"=25.11; "AE"=1080; "AF"=3.3; "AG"=30; "AH"=0.84; "AI"=19.51; "AJ"=59991641 }
$rowData[4] = @{ "D"=3941; "E"=24; "F"=24; "G"=24; "H"=19; "I"=19; "J"=0; "K"=1527; "L"=938; "M"=590; "N"=583; "O"=6; "P"=300; "Q"=-191; "R"=-58; "S"=176; "T"=43; "U"=-233; "V"=135; "W"=0.6; "X"=0.47; "Y"=3.17; "Z"=1.34; "AA"=158.98; "AB"=112.88; "AC"=31; "AD"=75.12; "AE"=1055; "AF"=2.22; "AG"=20; "AH"=0.85; "AI"=59.06; "AJ"=59991641 }
$rowData[5] = @{ "D"=4345; "E"=35; "F"=35; "G"=35; "H"=21; "I"=22; "J"=-1; "K"=1714; "L"=1112; "M"=603; "N"=597; "O"=5; "P"=300; "Q"=63; "R"=-31; "S"=45; "T"=50; "U"=13; "V"=135; "W"=0.8; "X"=0.47; "Y"=3.67; "Z"=1.27; "AA"=184.5; "AB"=118.43; "AC"=36; "AD"=51.86; "AE"=1080; "AF"=1.73; "AG"=30; "AH"=1.6; "AI"=76.69; "AJ"=59991641 }
$rowData[6] = @{ "D"=4457; "E"=9; "F"=9; "G"=17; "H"=12; "I"=14; "K"=1554; "L"=960; "M"=595; "N"=591; "P"=300; "Q"=-37; "R"=-44; "S"=23; "T"=51; "U"=-88; "V"=135; "W"=0.21; "X"=0.26; "Y"=2.28; "Z"=0.71; "AA"=161.38; "AB"=116.27; "AC"=23; "AD"=60.82; "AE"=1069; "AF"=1.29; "AG"=10; "AH"=0.73; "AI"=40.78; "AJ"=59991641 }

foreach ($r in $rowData.Keys) {
    $cols = $rowData[$r]
    foreach ($c in $cols.Keys) {
        $ws.Range($c + $r).Value = $cols[$c]
    }
}

# Remove erroneous estimate rows (2019E-2021E), keep identifying columns A-C only
$ws.Range("D7:AJ9").ClearContents()